# Weekly fruit/vegetable price update.
#
# A new daily record (Fecha 45275, i.e. 2023-12-15) is inserted as the new
# first data row for this product/market sheet. All the previously existing
# rows (172-206) shift down by one (173-207); the sheet therefore grows from
# 206 to 207 data-bearing rows (A1:R206 -> A1:R207).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row at row 172, pushing the old row 172 (and everything
# below it) down to row 173, etc. Excel copies the row-above's formatting
# (including the date number format on column D) into the new row.
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(172, 1).Value  = 11
$ws.Cells.Item(172, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(172, 3).Value  = "Bíobío"
$ws.Cells.Item(172, 4).Value  = 45275
$ws.Cells.Item(172, 5).Value  = 8
$ws.Cells.Item(172, 6).Value  = 100112001
$ws.Cells.Item(172, 7).Value  = "Berenjena"
$ws.Cells.Item(172, 8).Value  = "Sin especificar"
$ws.Cells.Item(172, 9).Value  = "Primera"
$ws.Cells.Item(172, 10).Value = 150
$ws.Cells.Item(172, 11).Value = 8000
$ws.Cells.Item(172, 12).Value = 8000
$ws.Cells.Item(172, 13).Value = 8000
$ws.Cells.Item(172, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(172, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(172, 16).Value = 133
$ws.Cells.Item(172, 17).Value = 60
$ws.Cells.Item(172, 18).Value = "Hortaliza"
